$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Footer "datetimeFigureOut" field: 04/06/2020 -> 05/06/2020
#    Lives on the Slide Master's Date Placeholder and on every one of
#    the 11 Custom Layouts' Date Placeholder shapes.
# ------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "04/06/2020") {
                $tr.Text = "05/06/2020"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ------------------------------------------------------------------
# 2) Slide 13 ("Confusion Matrix"): swap the FP / FN counts.
#    False Positive (FP) = 11  ->  False Positive (FP) = 51
#    False Negative (FN) = 51  ->  False Negative (FN) = 11
# ------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
for ($i = 1; $i -le $s13.Shapes.Count; $i++) {
    $shp = $s13.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -like "*False Positive (FP)*") {
        $tr = $shp.TextFrame.TextRange

        $full = $tr.Text
        $label = "False Positive (FP) = "
        $idx = $full.IndexOf($label)
        $numStart = $idx + $label.Length
        $numChars = $tr.Characters($numStart + 1, 3)
        $numChars.Text = "51 "

        $full = $tr.Text
        $label2 = "False Negative (FN) = "
        $idx2 = $full.IndexOf($label2)
        $numStart2 = $idx2 + $label2.Length
        $numChars2 = $tr.Characters($numStart2 + 1, 2)
        $numChars2.Text = "11"
    }
}

# ------------------------------------------------------------------
# 3) Slide 6: "Insert New Feature As Data's Prediction"
#              -> "Insert New Feature As Data's Class"
#    (TextBox 3, nested inside Group 1)
# ------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$grp6 = $s6.Shapes.Item(1)
$items6 = $grp6.GroupItems
for ($i = 1; $i -le $items6.Count; $i++) {
    $shp = $items6.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -like "*Insert New Feature*") {
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf("Prediction")
        $chars = $tr.Characters($idx + 1, $full.Length - $idx)
        $chars.Text = "Class"
    }
}

# ------------------------------------------------------------------
# 4) Slide 9: "Column Body as feature" -> "Column Body as feature / attribute"
# ------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
for ($i = 1; $i -le $s9.Shapes.Count; $i++) {
    $shp = $s9.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -like "*Column Body as feature*") {
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf("feature")
        $chars = $tr.Characters($idx + 1, 7)
        $chars.Text = "feature / attribute"
    }
}
